# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D) for the 69f88810-... row
# (row 5) on both the zh-cn and de-de status sheets to reflect the new
# handoff that just happened.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D5").Value = "2016-03-04 05:29:35"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D5").Value = "2016-03-04 05:29:52"
